$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.176.03"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "3.569.49"
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'606.33"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").Value = "'144.61"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("D7").Value = "3.569.00"
$ws.Range("E7").Value = "  +2.39%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.490"
$ws.Range("E9").Value = "  +3.03%  "

$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("E11").Value = "  -3.00%  "

$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("D13").Value = "4.174.10"
$ws.Range("E13").Value = "  +2.46%  "

$ws.Range("E14").Value = "  +2.33%  "

$ws.Range("D15").Value = "'30.16"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").Value = "3.585.42"
$ws.Range("E16").Value = "  +3.04%  "

$ws.Range("D17").Value = "66.256.28"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "'0.115"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").Value = "'11.48"
$ws.Range("E19").Value = "  +5.99%  "

$ws.Range("D20").Value = "'6.22"
$ws.Range("E20").Value = "  +1.27%  "

$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("D22").Value = "'431.38"
$ws.Range("E22").Value = "  +1.26%  "

$ws.Range("E23").Value = "  +2.68%  "

$ws.Range("D24").Value = "'79.60"
$ws.Range("E24").Value = "  +2.24%  "

$ws.Range("D25").Value = "3.712.82"
$ws.Range("E25").Value = "  +2.36%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("E27").Value = "  +1.59%  "

$ws.Range("E28").Value = "  +2.13%  "

$ws.Range("D29").Value = "'9.14"
$ws.Range("E29").Value = "  -1.21%  "

$ws.Range("D30").Value = "'7.90"
$ws.Range("E30").Value = "  -0.58%  "

$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("D32").Value = "3.566.55"
$ws.Range("E32").Value = "  +2.81%  "

$ws.Range("D33").Value = "'25.50"
$ws.Range("E33").Value = "  +1.74%  "

$ws.Range("D34").Value = "'1.45"
$ws.Range("E34").Value = "  -1.12%  "

$ws.Range("D35").Value = "'0.152"
$ws.Range("E35").Value = "  -8.03%  "

$ws.Range("E36").Value = "  +1.52%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  -0.51%  "

$ws.Range("E39").Value = "  -0.48%  "

$ws.Range("D40").Value = "'174.02"
$ws.Range("E40").Value = "  +2.21%  "

$ws.Range("D41").Value = "'0.0846"
$ws.Range("E41").Value = "  -1.52%  "

$ws.Range("D42").Value = "'5.20"
$ws.Range("E42").Value = "  +0.89%  "

$ws.Range("D43").Value = "'0.895"
$ws.Range("E43").Value = "  +1.80%  "

$ws.Range("E44").Value = "  +2.19%  "

$ws.Range("D45").Value = "'46.02"
$ws.Range("E45").Value = "  +1.35%  "

$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("D48").Value = "'25.07"
$ws.Range("E48").Value = "  -4.02%  "

$ws.Range("E49").Value = "  +0.93%  "

$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("D51").Value = "'23.13"
$ws.Range("E51").Value = "  +6.41%  "
